# Updated test data for German, Czech market. Added test data for Belgium market.
$wb = $excel.ActiveWorkbook

$wsGermany = $wb.Worksheets.Item("Germany")
$wsBelgium = $wb.Worksheets.Item("Belgium")
$wsCzech   = $wb.Worksheets.Item("Czech")

# Append the extra ticket references to the Belgium and Germany user-story cells.
# (Belgium is updated first so the new shared-string entries land in the same
# order as the source workbook: Belgium's string before Germany's.)
$wsBelgium.Range("B4").Value = "NGC-3478/T2265/T2267"
$wsGermany.Range("B4").Value = "NGC-3475/T1730/T1746"

# Move the active selection on the Germany and Belgium sheets to B5.
$wsGermany.Range("B5").Select()
$wsBelgium.Range("B5").Select()

# Germany becomes the active/visible sheet when the workbook is reopened
# (previously Czech was active).
$wsGermany.Activate()
